$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated Cypher queries (column B) ---

$qParticipants = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.library_strategy in ['Bisulfite-Seq']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id limit 100"

$qSamples = "MATCH (samp:sample)-->(p:participant)-->(s:study)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag`nWHERE g.library_strategy in ['Bisulfite-Seq']`nWITH DISTINCT s, p, samp`nRETURN`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(samp.sample_tumor_status,'') as ``Tumor``,`n    coalesce(samp.sample_type,'') as ``Analyte Type```nORDER BY samp.sample_id limit 100"

$qFiles = "MATCH (f:file)-->(s:study)`nOPTIONAL MATCH (samp:sample)<--(f)`nOPTIONAL MATCH (samp)-->(p:participant)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag`nWHERE g.library_strategy in ['Bisulfite-Seq']`nWITH DISTINCT f, s, p, samp`nRETURN`n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name,'') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id, '') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER BY f.file_name limit 100"

$qStat = "CALL{`n    MATCH (p:participant)-->(s:study)`n    OPTIONAL MATCH (samp:sample)-->(p)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag`n    WHERE g.library_strategy in ['Bisulfite-Seq']`n    RETURN `n        count(distinct p) AS num_participants`n}`nWITH num_participants`nCALL {`n    MATCH (samp:sample)-->(p:participant)-->(s)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag`n    WHERE g.library_strategy in ['Bisulfite-Seq']`n    RETURN `n        count(distinct samp) AS num_samples`n}`nWITH num_participants, num_samples`nCALL {`n    MATCH (f:file)-->(s:study)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (samp:sample)<--(f)`n    OPTIONAL MATCH (p:participant)<--(samp)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag`n    WHERE g.library_strategy in ['Bisulfite-Seq']`n    RETURN `n        count(distinct s) AS num_studies,`n        count(distinct f) AS num_files`n}`nRETURN `n    num_studies AS Studies,`n    num_participants AS Participants,`n    num_samples AS Samples,`n    num_files AS ``Files``"

# --- Row 2: was "CasesTab" -> now "ParticipantsTab", with updated queries ---
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $qParticipants
$ws.Range("C2").Value = $qStat

# --- Row 3: "SamplesTab" keeps its label, query text updated ---
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $qSamples
$ws.Range("C3").Value = $qStat

# --- Row 4: "FilesTab" keeps its label, query text updated ---
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $qFiles
$ws.Range("C4").Value = $qStat

# Writing the long wrapped-text queries makes the engine recompute
# row heights for the affected rows; restore the original (unchanged)
# custom row heights.
$ws.Rows.Item(2).RowHeight = 242.25
$ws.Rows.Item(3).RowHeight = 260.25
$ws.Rows.Item(4).RowHeight = 279.75

# Column A is best-fit / autosized to its text; "ParticipantsTab" is
# longer than the old "CasesTab" label, so the column widens to fit.
$ws.Columns.Item(1).AutoFit()

# --- Selection moves to a single cell D4 (was D4:E4) ---
$ws.Range("D4").Select()
